# FormulaTest.xlsx - "ValidateFormulas" sheet
# Insert a new row above the old row 31 ("Boolean" section) and fill it in
# with a small "If" demo block (mirrors the existing "Boolean" block style),
# then leave the selection on the newly added D31 formula cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidateFormulas")

# Shift row 31 (and everything below it) down by one row.
$ws.Rows("31:31").Insert()

# New header cell, styled like the other bold section headers (e.g. A29 "Vlookup").
$ws.Range("A31").Value = "If"
$ws.Range("A31").Font.Bold = $true

# New formulas demonstrating IF().
$ws.Range("B31").Formula = "=IF(B2>3,B3,B5)"
$ws.Range("C31").Formula = "=IF((B2*B3)*C1<0,(B2*B3)*C1,ABS((B2*B3)*C1))"
$ws.Range("D31").Formula = "=IF((B2*B3)*C1<0,ABS((B2*B3)*C1),(B2*B3)*C1)"

# Match the recorded selection after the edit.
$ws.Range("D31").Select()
